$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.5725316666666668
$ws.Range("H2").Value = 1.717595
$ws.Range("I2").Value = 0.3864899584549088
$ws.Range("J2").Value = 0.3864899584549088
$ws.Range("M2").Value = 16.28844733333333
$ws.Range("N2").Value = 48.865342
$ws.Range("O2").Value = 0.2176904746803693
$ws.Range("P2").Value = 0.2176904746803693
$ws.Range("Q2").Value = 9.325651899165557
$ws.Range("R2").Value = 83.93086709249
$ws.Range("S2").Value = 0.08413518251524532
$ws.Range("T2").Value = 0.08413518251524532
$ws.Range("G3").Value = 0.5725316666666668
$ws.Range("H3").Value = 1.717595
$ws.Range("I3").Value = 0.3864899584549088
$ws.Range("J3").Value = 0.3864899584549088
$ws.Range("M3").Value = 27.61090666666666
$ws.Range("N3").Value = 82.83271999999999
$ws.Range("O3").Value = 0.3690119294748028
$ws.Range("P3").Value = 0.3690119294748029
$ws.Range("Q3").Value = 15.80811841204445
$ws.Range("R3").Value = 142.2730657084
$ws.Range("S3").Value = 0.1426194052920823
$ws.Range("T3").Value = 0.1426194052920823
$ws.Range("G4").Value = 0.5725316666666668
$ws.Range("H4").Value = 1.717595
$ws.Range("I4").Value = 0.3864899584549088
$ws.Range("J4").Value = 0.3864899584549088
$ws.Range("M4").Value = 26.266325
$ws.Range("N4").Value = 78.798975
$ws.Range("O4").Value = 0.3510419771967738
$ws.Range("P4").Value = 0.3510419771967739
$ws.Range("Q4").Value = 15.03830282945833
$ws.Range("R4").Value = 135.344725465125
$ws.Range("S4").Value = 0.1356741991827101
$ws.Range("T4").Value = 0.1356741991827102
$ws.Range("G5").Value = 0.5725316666666668
$ws.Range("H5").Value = 1.717595
$ws.Range("I5").Value = 0.3864899584549088
$ws.Range("J5").Value = 0.3864899584549088
$ws.Range("M5").Value = 4.658207333333333
$ws.Range("N5").Value = 13.974622
$ws.Range("O5").Value = 0.06225561864805391
$ws.Range("P5").Value = 0.06225561864805392
$ws.Range("Q5").Value = 2.666971208232223
$ws.Range("R5").Value = 24.00274087409
$ws.Range("S5").Value = 0.024061171464871
$ws.Range("T5").Value = 0.02406117146487101
$ws.Range("H6").Value = 0.919331
$ws.Range("I6").Value = 0.2068661122070742
$ws.Range("J6").Value = 0.2068661122070743
$ws.Range("M6").Value = 16.28844733333333
$ws.Range("N6").Value = 48.865342
$ws.Range("O6").Value = 0.2176904746803693
$ws.Range("P6").Value = 0.2176904746803693
$ws.Range("Q6").Value = 4.991491525133556
$ws.Range("R6").Value = 44.923423726202
$ws.Range("S6").Value = 0.04503278216164053
$ws.Range("T6").Value = 0.04503278216164054
$ws.Range("H7").Value = 0.919331
$ws.Range("I7").Value = 0.2068661122070742
$ws.Range("J7").Value = 0.2068661122070743
$ws.Range("M7").Value = 27.61090666666666
$ws.Range("N7").Value = 82.83271999999999
$ws.Range("O7").Value = 0.3690119294748028
$ws.Range("P7").Value = 0.3690119294748029
$ws.Range("Q7").Value = 8.461187478924444
$ws.Range("R7").Value = 76.15068731032
$ws.Range("S7").Value = 0.07633606320848353
$ws.Range("T7").Value = 0.07633606320848355
$ws.Range("H8").Value = 0.919331
$ws.Range("I8").Value = 0.2068661122070742
$ws.Range("J8").Value = 0.2068661122070743
$ws.Range("M8").Value = 26.266325
$ws.Range("N8").Value = 78.798975
$ws.Range("O8").Value = 0.3510419771967738
$ws.Range("P8").Value = 0.3510419771967739
$ws.Range("Q8").Value = 8.049148942858332
$ws.Range("R8").Value = 72.442340485725
$ws.Range("S8").Value = 0.07261868904418101
$ws.Range("T8").Value = 0.07261868904418103
$ws.Range("H9").Value = 0.919331
$ws.Range("I9").Value = 0.2068661122070742
$ws.Range("J9").Value = 0.2068661122070743
$ws.Range("M9").Value = 4.658207333333333
$ws.Range("N9").Value = 13.974622
$ws.Range("O9").Value = 0.06225561864805391
$ws.Range("P9").Value = 0.06225561864805392
$ws.Range("Q9").Value = 1.427478135320222
$ws.Range("R9").Value = 12.847303217882
$ws.Range("S9").Value = 0.01287857779276914
$ws.Range("T9").Value = 0.01287857779276915
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.5451493333333333
$ws.Range("H10").Value = 1.635448
$ws.Range("I10").Value = 0.3680053968340403
$ws.Range("J10").Value = 0.3680053968340404
$ws.Range("M10").Value = 16.28844733333333
$ws.Range("N10").Value = 48.865342
$ws.Range("O10").Value = 0.2176904746803693
$ws.Range("P10").Value = 0.2176904746803693
$ws.Range("Q10").Value = 8.879636204801777
$ws.Range("R10").Value = 79.91672584321599
$ws.Range("S10").Value = 0.08011126952173993
$ws.Range("T10").Value = 0.08011126952173994
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.5451493333333333
$ws.Range("H11").Value = 1.635448
$ws.Range("I11").Value = 0.3680053968340403
$ws.Range("J11").Value = 0.3680053968340404
$ws.Range("M11").Value = 27.61090666666666
$ws.Range("N11").Value = 82.83271999999999
$ws.Range("O11").Value = 0.3690119294748028
$ws.Range("P11").Value = 0.3690119294748029
$ws.Range("Q11").Value = 15.05206736206222
$ws.Range("R11").Value = 135.46860625856
$ws.Range("S11").Value = 0.1357983815428697
$ws.Range("T11").Value = 0.1357983815428697
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.5451493333333333
$ws.Range("H12").Value = 1.635448
$ws.Range("I12").Value = 0.3680053968340403
$ws.Range("J12").Value = 0.3680053968340404
$ws.Range("M12").Value = 26.266325
$ws.Range("N12").Value = 78.798975
$ws.Range("O12").Value = 0.3510419771967738
$ws.Range("P12").Value = 0.3510419771967739
$ws.Range("Q12").Value = 14.31906956286666
$ws.Range("R12").Value = 128.8716260658
$ws.Range("S12").Value = 0.1291853421237049
$ws.Range("T12").Value = 0.1291853421237049
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.5451493333333333
$ws.Range("H13").Value = 1.635448
$ws.Range("I13").Value = 0.3680053968340403
$ws.Range("J13").Value = 0.3680053968340404
$ws.Range("M13").Value = 4.658207333333333
$ws.Range("N13").Value = 13.974622
$ws.Range("O13").Value = 0.06225561864805391
$ws.Range("P13").Value = 0.06225561864805392
$ws.Range("Q13").Value = 2.539418622295111
$ws.Range("R13").Value = 22.854767600656
$ws.Range("S13").Value = 0.02291040364572576
$ws.Range("T13").Value = 0.02291040364572577
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.05723766666666667
$ws.Range("H14").Value = 0.171713
$ws.Range("I14").Value = 0.03863853250397663
$ws.Range("J14").Value = 0.03863853250397663
$ws.Range("M14").Value = 16.28844733333333
$ws.Range("N14").Value = 48.865342
$ws.Range("O14").Value = 0.2176904746803693
$ws.Range("P14").Value = 0.2176904746803693
$ws.Range("Q14").Value = 0.9323127189828889
$ws.Range("R14").Value = 8.390814470845999
$ws.Range("S14").Value = 0.008411240481743552
$ws.Range("T14").Value = 0.008411240481743553
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.05723766666666667
$ws.Range("H15").Value = 0.171713
$ws.Range("I15").Value = 0.03863853250397663
$ws.Range("J15").Value = 0.03863853250397663
$ws.Range("M15").Value = 27.61090666666666
$ws.Range("N15").Value = 82.83271999999999
$ws.Range("O15").Value = 0.3690119294748028
$ws.Range("P15").Value = 0.3690119294748029
$ws.Range("Q15").Value = 1.580383872151111
$ws.Range("R15").Value = 14.22345484936
$ws.Range("S15").Value = 0.0142580794313673
$ws.Range("T15").Value = 0.0142580794313673
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.05723766666666667
$ws.Range("H16").Value = 0.171713
$ws.Range("I16").Value = 0.03863853250397663
$ws.Range("J16").Value = 0.03863853250397663
$ws.Range("M16").Value = 26.266325
$ws.Range("N16").Value = 78.798975
$ws.Range("O16").Value = 0.3510419771967738
$ws.Range("P16").Value = 0.3510419771967739
$ws.Range("Q16").Value = 1.503423154908333
$ws.Range("R16").Value = 13.530808394175
$ws.Range("S16").Value = 0.01356374684617777
$ws.Range("T16").Value = 0.01356374684617777
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.05723766666666667
$ws.Range("H17").Value = 0.171713
$ws.Range("I17").Value = 0.03863853250397663
$ws.Range("J17").Value = 0.03863853250397663
$ws.Range("M17").Value = 4.658207333333333
$ws.Range("N17").Value = 13.974622
$ws.Range("O17").Value = 0.06225561864805391
$ws.Range("P17").Value = 0.06225561864805392
$ws.Range("Q17").Value = 0.2666249186095556
$ws.Range("R17").Value = 2.399624267486
$ws.Range("S17").Value = 0.002405465744688005
$ws.Range("T17").Value = 0.002405465744688006

Write-Output "Applied 204 cell updates"
